$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.887.74'
$ws.Range("E2").Value = '  +0.72%  '

$ws.Range("D3").Value = '1.627.89'
$ws.Range("E3").Value = '  +1.00%  '

$ws.Range("D5").Value = '''214.40'
$ws.Range("E5").Value = '  +0.91%  '

$ws.Range("E6").Value = '  -0.77%  '

$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("D8").Value = '''28.45'
$ws.Range("E8").Value = '  -1.25%  '

$ws.Range("D9").Value = '''0.257'
$ws.Range("E9").Value = '  -0.20%  '

$ws.Range("D10").Value = '''0.0607'
$ws.Range("E10").Value = '  -0.11%  '

$ws.Range("D11").Value = '''0.0901'
$ws.Range("E11").Value = '  -0.53%  '

$ws.Range("D12").Value = '1.863.77'
$ws.Range("E12").Value = '  +0.95%  '

$ws.Range("D13").Value = '1.635.40'
$ws.Range("E13").Value = '  +1.08%  '

$ws.Range("D14").Value = '''0.560'
$ws.Range("E14").Value = '  -0.62%  '

$ws.Range("D15").Value = '''9.24'
$ws.Range("E15").Value = '  +7.07%  '

$ws.Range("D16").Value = '29.907.62'
$ws.Range("E16").Value = '  +0.62%  '

$ws.Range("D17").Value = '''3.82'
$ws.Range("E17").Value = '  -0.59%  '

$ws.Range("D18").Value = '''63.85'
$ws.Range("E18").Value = '  -1.29%  '

$ws.Range("D19").Value = '''240.35'
$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").Value = '0.0₃0699'
$ws.Range("E20").Value = '  -0.73%  '

$ws.Range("E21").Value = '  +0.54%  '

$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '''9.78'
$ws.Range("E22").Value = '  +2.05%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '''4.11'
$ws.Range("E23").Value = '  +0.69%  '

$ws.Range("D24").Value = '''2.16'
$ws.Range("E24").Value = '  +1.37%  '

$ws.Range("D25").Value = '''157.76'
$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("D26").Value = '''15.41'
$ws.Range("E26").Value = '  -1.11%  '

$ws.Range("E27").Value = '  -0.64%  '

$ws.Range("D28").Value = '''6.53'
$ws.Range("E28").Value = '  -0.27%  '

$ws.Range("E29").Value = '  +0.46%  '

$ws.Range("D30").Value = '''0.0486'
$ws.Range("E30").Value = '  +1.06%  '

$ws.Range("E31").Value = '  +3.16%  '

$ws.Range("D32").Value = '''3.36'
$ws.Range("E32").Value = '  +1.67%  '

$ws.Range("D33").Value = '''3.16'
$ws.Range("E33").Value = '  -0.53%  '

$ws.Range("D34").Value = '1.423.74'
$ws.Range("E34").Value = '  -2.14%  '

$ws.Range("E35").Value = '  +3.88%  '

$ws.Range("E36").Value = '  -1.77%  '

$ws.Range("D37").Value = '''2.74'
$ws.Range("E37").Value = '  -5.17%  '

$ws.Range("E38").Value = '  +0.00%  '

$ws.Range("E39").Value = '  -0.16%  '

$ws.Range("D40").Value = '''74.69'
$ws.Range("E40").Value = '  +7.35%  '

$ws.Range("D41").Value = '''0.553'
$ws.Range("E41").Value = '  -0.48%  '

$ws.Range("E42").Value = '  -1.28%  '

$ws.Range("D43").Value = '''0.0497'
$ws.Range("E43").Value = '  -1.67%  '

$ws.Range("D44").Value = '''0.825'
$ws.Range("E44").Value = '  +0.22%  '

$ws.Range("E45").Value = '  +0.69%  '

$ws.Range("E46").Value = '  +0.31%  '

$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.771.28'
$ws.Range("E47").Value = '  +0.93%  '

$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").Value = '''49.60'
$ws.Range("E48").Value = '  -7.04%  '

$ws.Range("E49").Value = '  -2.40%  '

$ws.Range("D50").Value = '''90.62'
$ws.Range("E50").Value = '  +3.89%  '

$ws.Range("E51").Value = '  +9.34%  '
